$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Формулировка задания для Варианта 16" ->
#    "Формулировка задания для Варианта " + "25" (own run) + moved
#    "_GoBack" bookmark right after the new "25" run.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Формулировка задания для Варианта 16", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$fullMatch = $rng1.Duplicate
$numberRange = $fullMatch.Duplicate
$numberRange.Start = $fullMatch.End - 2
$numberRange.End = $fullMatch.End

# Remove the old "16" and type the new "25" in its place.
$numberRange.Delete()
$insertionPoint = $numberRange.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertAfter("25")

# Force the newly typed "25" to live in its own run (it currently has
# identical formatting to its neighbour and would otherwise be silently
# merged back together) by toggling an unrelated character property.
$newNumberRange = $d.Range($insertionPoint.Start, $insertionPoint.Start + 2)
$newNumberRange.Font.Bold = $true
$newNumberRange.Font.Bold = $false

# Word keeps only a single "_GoBack" bookmark - re-adding it here moves
# it from wherever it used to be (see step 3 below).
$bmPoint = $d.Range($newNumberRange.End, $newNumberRange.End)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ---------------------------------------------------------------------
# 2) "Рисунок 4 - Результат выполнения" caption: drop the stray
#    <w:lang w:val="en-US"/> override that sits on the "4" run.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Рисунок 4 - Результат выполнения", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$digitRange = $d.Range($rng2.Start + 8, $rng2.Start + 9)
$digitRange.Delete()

# Retype the digit right after the preceding "Рисунок " run so it picks
# up that run's (language-less) formatting instead of the old one.
$digitInsertionPoint = $d.Range($rng2.Start + 8, $rng2.Start + 8)
$digitInsertionPoint.InsertAfter("4")

# It would now be silently merged into "Рисунок " (identical formatting);
# split it back into its own run, matching the original structure.
$newDigitRange = $d.Range($rng2.Start + 8, $rng2.Start + 9)
$newDigitRange.Font.Bold = $true
$newDigitRange.Font.Bold = $false

# ---------------------------------------------------------------------
# 3) The old "_GoBack" bookmark (next to the "Рисунок 5" caption) is
#    implicitly removed by step 1 re-adding the bookmark elsewhere,
#    since Word only ever keeps one "_GoBack" bookmark at a time.
# ---------------------------------------------------------------------
